$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 373, shifting existing rows 373:397 down to 374:398
$ws.Rows.Item(373).Insert()

# Populate the new row 373 with the new weekly price entry
$ws.Cells.Item(373, 1).Value = 11
$ws.Cells.Item(373, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(373, 3).Value = "Bíobío"
$ws.Cells.Item(373, 4).Value = 45013
$ws.Cells.Item(373, 5).Value = 8
$ws.Cells.Item(373, 6).Value = 100112008
$ws.Cells.Item(373, 7).Value = "Coliflor"
$ws.Cells.Item(373, 8).Value = "Sin especificar"
$ws.Cells.Item(373, 9).Value = "Primera"
$ws.Cells.Item(373, 10).Value = 1000
$ws.Cells.Item(373, 11).Value = 1500
$ws.Cells.Item(373, 12).Value = 1500
$ws.Cells.Item(373, 13).Value = 1500
$ws.Cells.Item(373, 14).Value = "`$/unidad"
$ws.Cells.Item(373, 15).Value = "Región Metropolitana"
$ws.Cells.Item(373, 16).Value = 1500
$ws.Cells.Item(373, 17).Value = 1
$ws.Cells.Item(373, 18).Value = "Hortaliza"
